# "Thailand on background of flow"
# Update the "loads" sheet (sheet3 / power-factor calc) with new flow-based
# P/Q/S figures, push the old power-factor-derived figures further down the
# sheet as a labelled reference block, and add a new labelled block with the
# new (flow-based) figures repeated for comparison.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("loads")

# --- 1. Preserve the old (pre-edit) power-factor-derived numbers for the
#        "Ruoya" reference block (rows 28:31) before they get overwritten in
#        rows 2:5. These keep the original bold/body style (s="1").
$ws.Range("D28:F28").Font.Name = "Aptos Narrow"
$ws.Range("D28").Value = 2.19
$ws.Range("E28").Value = 1.0606654099999999
$ws.Range("F28").Value = 2.43333333

$ws.Range("D29:F29").Font.Name = "Aptos Narrow"
$ws.Range("D29").Value = 2.99
$ws.Range("E29").Value = 1.44812309
$ws.Range("F29").Value = 3.32222222

$ws.Range("D30:F30").Font.Name = "Aptos Narrow"
$ws.Range("D30").Value = 2.21
$ws.Range("E30").Value = 1.07035185
$ws.Range("F30").Value = 2.4555555600000001

$ws.Range("D31:F31").Font.Name = "Aptos Narrow"
$ws.Range("D31").Value = 19.83
$ws.Range("E31").Value = 9.60410734
$ws.Range("F31").Value = 22.0333333

# Label above the old block.
$ws.Range("D27").Value = "Ruoya"

# --- 2. New "Willi" labelled block (rows 21:24) mirrors the new flow-based
#        numbers that are about to replace D2:F5, but with default/no style.
$ws.Range("D21").Value = 4.4
$ws.Range("E21").Value = 2.1
$ws.Range("F21").Value = 4.9

$ws.Range("D22").Value = 6.1
$ws.Range("E22").Value = 2.9
$ws.Range("F22").Value = 6.8

$ws.Range("D23").Value = 4.1
$ws.Range("E23").Value = 1.9
$ws.Range("F23").Value = 4.5

$ws.Range("D24").Value = 4.5
$ws.Range("E24").Value = 2.2
$ws.Range("F24").Value = 4.9

# Label above the new block.
$ws.Range("D20").Value = "Willi"

# --- 3. New aggregate flow totals (rows 16:19), keeping the body style.
$ws.Range("D16:F16").Font.Name = "Aptos Narrow"
$ws.Range("D16").Value = 4422.44898
$ws.Range("E16").Value = 2141.8898
$ws.Range("F16").Value = 4913.26413

$ws.Range("D17:F17").Font.Name = "Aptos Narrow"
$ws.Range("D17").Value = 6041.58163
$ws.Range("E17").Value = 2926.07153
$ws.Range("F17").Value = 6712.09243

$ws.Range("D18:F18").Font.Name = "Aptos Narrow"
$ws.Range("D18").Value = 40077.0408
$ws.Range("E18").Value = 19410.1968
$ws.Range("F18").Value = 44524.8974

$ws.Range("D19:F19").Font.Name = "Aptos Narrow"
$ws.Range("D19").Value = 4461.73469
$ws.Range("E19").Value = 2160.91674
$ws.Range("F19").Value = 4956.90987

# --- 4. Replace the headline D2:F5 figures with the new flow-based values
#        and strip their formatting back to the default/no style.
$ws.Range("D2:F5").ClearFormats()

$ws.Range("D2").Value = 4.4
$ws.Range("E2").Value = 2.1
$ws.Range("F2").Value = 4.9

$ws.Range("D3").Value = 6.1
$ws.Range("E3").Value = 2.9
$ws.Range("F3").Value = 6.8

$ws.Range("D4").Value = 4.1
$ws.Range("E4").Value = 1.9
$ws.Range("F4").Value = 4.5

$ws.Range("D5").Value = 4.5
$ws.Range("E5").Value = 2.2
$ws.Range("F5").Value = 4.9

# --- 5. Match the author's final view/selection state.
$ws.Application.ActiveWindow.ScrollRow = 6
$ws.Range("D21").Select()
